# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-25 (row 22 unchanged)
$kValues = @{
    2  = 4
    3  = 2
    4  = 8
    5  = 1
    6  = 6
    7  = 6
    8  = 4
    9  = 3
    10 = 4
    11 = 11
    12 = 9
    13 = 4
    14 = 0
    15 = 3
    16 = 4
    17 = 5
    18 = 5
    19 = 2
    20 = 4
    21 = 2
    23 = 1
    24 = 1
    25 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
